$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Update the "Vreme skeniranja" (scan time) column E for rows 2-12
# with the new timestamps recorded after finishing the scan.
$ws.Range("E2").Value  = "2025-04-10 18:27:32"
$ws.Range("E3").Value  = "2025-04-10 18:27:32"
$ws.Range("E4").Value  = "2025-04-10 18:27:34"
$ws.Range("E5").Value  = "2025-04-10 18:27:29"
$ws.Range("E6").Value  = "2025-04-10 18:27:29"
$ws.Range("E7").Value  = "2025-04-10 18:27:33"
$ws.Range("E8").Value  = "2025-04-10 18:27:28"
$ws.Range("E9").Value  = "2025-04-10 18:27:27"
$ws.Range("E10").Value = "2025-04-10 18:27:34"
$ws.Range("E11").Value = "2025-04-10 18:27:35"
$ws.Range("E12").Value = "2025-04-10 18:27:36"

# Remove the trailing rows (13-15) that belonged to boxes not found
# in the database - they are no longer part of the finished report.
$ws.Rows.Item(13).Delete()
$ws.Rows.Item(13).Delete()
$ws.Rows.Item(13).Delete()
